$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($cell, $value) {
    $cell.Value = $value
}

function Set-TextNumber($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-PlainText $ws.Range('D2') '62.753.56'
Set-PlainText $ws.Range('E2') '  +5.83%  '
Set-PlainText $ws.Range('D3') '3.475.60'
Set-PlainText $ws.Range('E3') '  +5.15%  '
Set-PlainText $ws.Range('E4') '  +0.12%  '
Set-TextNumber $ws.Range('D5') '408.84'
Set-PlainText $ws.Range('E5') '  +0.37%  '
Set-TextNumber $ws.Range('D6') '129.06'
Set-PlainText $ws.Range('E6') '  +15.90%  '
Set-PlainText $ws.Range('D7') '3.467.89'
Set-PlainText $ws.Range('E7') '  +5.13%  '
Set-TextNumber $ws.Range('D8') '0.596'
Set-PlainText $ws.Range('E8') '  +2.47%  '
Set-PlainText $ws.Range('E9') '  +0.11%  '
Set-TextNumber $ws.Range('D10') '0.690'
Set-PlainText $ws.Range('E10') '  +9.86%  '
Set-TextNumber $ws.Range('D11') '0.130'
Set-PlainText $ws.Range('E11') '  +32.95%  '
Set-TextNumber $ws.Range('D12') '42.79'
Set-PlainText $ws.Range('E12') '  +8.02%  '
Set-PlainText $ws.Range('D13') '4.028.63'
Set-PlainText $ws.Range('E13') '  +5.33%  '
Set-PlainText $ws.Range('E14') '  -0.67%  '
Set-TextNumber $ws.Range('D15') '8.71'
Set-PlainText $ws.Range('E15') '  +2.96%  '
Set-TextNumber $ws.Range('D16') '20.02'
Set-PlainText $ws.Range('E16') '  +3.98%  '
Set-PlainText $ws.Range('D17') '3.483.26'
Set-PlainText $ws.Range('E17') '  +5.22%  '
Set-PlainText $ws.Range('D18') '62.700.70'
Set-PlainText $ws.Range('E18') '  +6.28%  '
Set-TextNumber $ws.Range('D19') '1.04'
Set-PlainText $ws.Range('E19') '  +0.98%  '
Set-TextNumber $ws.Range('D20') '10.88'
Set-PlainText $ws.Range('E20') '  +2.11%  '
Set-PlainText $ws.Range('E21') '  +23.78%  '
Set-PlainText $ws.Range('E22') '  +0.48%  '
Set-TextNumber $ws.Range('D23') '81.83'
Set-PlainText $ws.Range('E23') '  +8.77%  '
Set-PlainText $ws.Range('E24') '  +0.04%  '
Set-TextNumber $ws.Range('D25') '309.88'
Set-PlainText $ws.Range('E25') '  +1.24%  '
Set-TextNumber $ws.Range('D26') '3.17'
Set-PlainText $ws.Range('E26') '  +0.17%  '
Set-TextNumber $ws.Range('D27') '30.31'
Set-PlainText $ws.Range('E27') '  +6.64%  '
Set-TextNumber $ws.Range('D28') '8.17'
Set-PlainText $ws.Range('E28') '  +5.41%  '
Set-TextNumber $ws.Range('D29') '7.73'
Set-PlainText $ws.Range('E29') '  +4.38%  '
Set-TextNumber $ws.Range('D30') '0.179'
Set-PlainText $ws.Range('E30') '  -0.91%  '
Set-TextNumber $ws.Range('D31') '4.38'
Set-PlainText $ws.Range('E31') '  -2.06%  '
Set-TextNumber $ws.Range('D32') '0.120'
Set-PlainText $ws.Range('E32') '  +4.40%  '
Set-PlainText $ws.Range('B33') 'Cosmos'
Set-PlainText $ws.Range('C33') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextNumber $ws.Range('D33') '11.93'
Set-PlainText $ws.Range('E33') '  +4.34%  '
Set-PlainText $ws.Range('B34') 'Toncoin'
Set-PlainText $ws.Range('C34') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextNumber $ws.Range('D34') '2.64'
Set-PlainText $ws.Range('E34') '  +24.10%  '
Set-TextNumber $ws.Range('D35') '42.85'
Set-PlainText $ws.Range('E35') '  +7.94%  '
Set-PlainText $ws.Range('E36') '  +0.05%  '
Set-TextNumber $ws.Range('D37') '0.0491'
Set-PlainText $ws.Range('E37') '  -3.30%  '
Set-TextNumber $ws.Range('D38') '52.54'
Set-PlainText $ws.Range('E38') '  +1.48%  '
Set-PlainText $ws.Range('E39') '  +5.97%  '
Set-TextNumber $ws.Range('D40') '0.997'
Set-PlainText $ws.Range('E40') '  -0.09%  '
Set-TextNumber $ws.Range('D41') '2.98'
Set-PlainText $ws.Range('E41') '  -5.15%  '
Set-PlainText $ws.Range('B42') 'Stellar'
Set-PlainText $ws.Range('C42') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextNumber $ws.Range('D42') '0.126'
Set-PlainText $ws.Range('E42') '  +2.91%  '
Set-PlainText $ws.Range('B43') 'ARBITRUM'
Set-PlainText $ws.Range('C43') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextNumber $ws.Range('D43') '1.98'
Set-PlainText $ws.Range('E43') '  +4.07%  '
Set-PlainText $ws.Range('B44') 'Monero'
Set-PlainText $ws.Range('C44') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextNumber $ws.Range('D44') '137.30'
Set-PlainText $ws.Range('E44') '  -0.99%  '
Set-TextNumber $ws.Range('D45') '17.42'
Set-PlainText $ws.Range('E45') '  +3.74%  '
Set-TextNumber $ws.Range('D46') '0.286'
Set-PlainText $ws.Range('E46') '  +2.03%  '
Set-TextNumber $ws.Range('D47') '3.93'
Set-PlainText $ws.Range('E47') '  +0.42%  '
Set-TextNumber $ws.Range('D48') '2.25'
Set-PlainText $ws.Range('E48') '  -0.76%  '
Set-TextNumber $ws.Range('D49') '22.33'
Set-PlainText $ws.Range('E49') '  +0.25%  '
Set-PlainText $ws.Range('D50') '2.206.19'
Set-PlainText $ws.Range('E50') '  +0.25%  '
Set-PlainText $ws.Range('B51') 'ApeXProtocol'
Set-PlainText $ws.Range('C51') 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextNumber $ws.Range('D51') '2.40'
Set-PlainText $ws.Range('E51') '  +0.07%  '

Write-Output "Applied all updates"
